$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "58.481.67"
$ws.Range("E2").Value = "  -2.87%  "
$ws.Range("D3").Value = "2.278.11"
$ws.Range("E3").Value = "  -5.97%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'545.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.36%  "
$ws.Range("D6").Value = "'130.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.83%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.569"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.51%  "
$ws.Range("E9").Value = "  -4.05%  "
$ws.Range("D10").Value = "'5.51"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.94%  "
$ws.Range("D11").Value = "'0.148"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("E12").Value = "  -5.70%  "
$ws.Range("E13").Value = "  -5.60%  "
$ws.Range("D14").Value = "2.683.68"
$ws.Range("E14").Value = "  -5.98%  "
$ws.Range("D15").Value = "58.414.80"
$ws.Range("E15").Value = "  -2.82%  "
$ws.Range("E16").Value = "  -3.69%  "
$ws.Range("D17").Value = "2.276.44"
$ws.Range("E17").Value = "  -6.33%  "
$ws.Range("D18").Value = "'10.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.25%  "
$ws.Range("D19").Value = "'4.29"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.64%  "
$ws.Range("D20").Value = "'314.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.37%  "
$ws.Range("E21").Value = "  -4.99%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "'62.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.22%  "
$ws.Range("E24").Value = "  -4.39%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("E26").Value = "  -6.90%  "
$ws.Range("E27").Value = "  -6.32%  "
$ws.Range("E28").Value = "  -1.25%  "
$ws.Range("D29").Value = "'170.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("D30").Value = "0.0₃0721"
$ws.Range("D31").Value = "'1.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("E32").Value = "  -5.76%  "
$ws.Range("D33").Value = "'0.383"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.55%  "
$ws.Range("D35").Value = "'17.79"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.42%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  -5.52%  "
$ws.Range("E38").Value = "  -6.51%  "
$ws.Range("D39").Value = "'38.13"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.01%  "
$ws.Range("D40").Value = "'1.52"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.65%  "
$ws.Range("D41").Value = "'303.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.66%  "
$ws.Range("D42").Value = "'140.78"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.64%  "
$ws.Range("D43").Value = "'3.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.50%  "
$ws.Range("D44").Value = "'0.0945"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.34%  "
$ws.Range("E45").Value = "  -3.70%  "
$ws.Range("D46").Value = "'0.550"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.88%  "
$ws.Range("D47").Value = "'18.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -9.03%  "
$ws.Range("D48").Value = "'0.0214"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.38%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'16.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.08%  "
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").Value = "'11.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("E51").Value = "  -0.61%  "
